# Apply the "Updated cryptos list" data refresh (Price / Volume(1h) columns,
# plus the RocketPoolETH/RenderToken row swap at the bottom of the table).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A leading apostrophe forces Excel to store a numeric-looking literal as
# text, matching the source workbook's Text-typed cells in the Price column.
$apos = "'"

$ws.Range('D2').Value = '36.491.51'
$ws.Range('E2').Value = '  -1.66%  '
$ws.Range('D3').Value = '2.057.00'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').Value = $apos + '242.25'
$ws.Range('E5').Value = '  -2.18%  '
$ws.Range('D6').Value = $apos + '0.666'
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = $apos + '54.67'
$ws.Range('E8').Value = '  -5.49%  '
$ws.Range('D9').Value = $apos + '58.48'
$ws.Range('E9').Value = '  -3.16%  '
$ws.Range('D10').Value = $apos + '0.358'
$ws.Range('E10').Value = '  -6.32%  '
$ws.Range('D11').Value = $apos + '0.0751'
$ws.Range('E11').Value = '  -3.72%  '
$ws.Range('E12').Value = '  -2.81%  '
$ws.Range('D13').Value = $apos + '0.892'
$ws.Range('E13').Value = '  -1.74%  '
$ws.Range('D14').Value = $apos + '14.71'
$ws.Range('E14').Value = '  -6.91%  '
$ws.Range('D15').Value = '2.355.92'
$ws.Range('E15').Value = '  -0.41%  '
$ws.Range('D16').Value = $apos + '5.35'
$ws.Range('E16').Value = '  -7.77%  '
$ws.Range('D17').Value = '2.036.07'
$ws.Range('E17').Value = '  -0.81%  '
$ws.Range('D18').Value = '36.429.98'
$ws.Range('E18').Value = '  -1.94%  '
$ws.Range('D19').Value = $apos + '16.77'
$ws.Range('E19').Value = '  -9.12%  '
$ws.Range('D20').Value = $apos + '72.04'
$ws.Range('E20').Value = '  -3.75%  '
$ws.Range('D21').Value = '0.0₃0856'
$ws.Range('E21').Value = '  -5.25%  '
$ws.Range('D22').Value = $apos + '238.31'
$ws.Range('E22').Value = '  +0.41%  '
$ws.Range('D23').Value = $apos + '5.25'
$ws.Range('E23').Value = '  -3.89%  '
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('E25').Value = '  -5.04%  '
$ws.Range('D26').Value = $apos + '9.37'
$ws.Range('E26').Value = '  -2.31%  '
$ws.Range('D27').Value = $apos + '2.12'
$ws.Range('E27').Value = '  -3.05%  '
$ws.Range('D28').Value = $apos + '162.84'
$ws.Range('E28').Value = '  -4.25%  '
$ws.Range('D29').Value = $apos + '20.16'
$ws.Range('E29').Value = '  +0.39%  '
$ws.Range('E30').Value = '  -2.21%  '
$ws.Range('D31').Value = $apos + '5.04'
$ws.Range('E31').Value = '  -7.44%  '
$ws.Range('D32').Value = $apos + '1.17'
$ws.Range('E32').Value = '  +2.90%  '
$ws.Range('D33').Value = $apos + '4.49'
$ws.Range('E33').Value = '  -6.59%  '
$ws.Range('D34').Value = $apos + '0.0592'
$ws.Range('E34').Value = '  -4.82%  '
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('D36').Value = $apos + '1.85'
$ws.Range('E36').Value = '  +1.24%  '
$ws.Range('D37').Value = $apos + '0.0824'
$ws.Range('E37').Value = '  -6.18%  '
$ws.Range('D38').Value = $apos + '2.18'
$ws.Range('E38').Value = '  -5.32%  '
$ws.Range('D39').Value = $apos + '1.24'
$ws.Range('E39').Value = '  -7.02%  '
$ws.Range('D40').Value = $apos + '4.86'
$ws.Range('E40').Value = '  -6.20%  '
$ws.Range('D41').Value = $apos + '0.0214'
$ws.Range('E41').Value = '  -4.49%  '
$ws.Range('E42').Value = '  -8.80%  '
$ws.Range('D43').Value = $apos + '1.11'
$ws.Range('E43').Value = '  -4.39%  '
$ws.Range('D44').Value = $apos + '93.59'
$ws.Range('E44').Value = '  -5.90%  '
$ws.Range('D45').Value = $apos + '0.0901'
$ws.Range('E45').Value = '  -10.48%  '
$ws.Range('D46').Value = '1.389.17'
$ws.Range('E46').Value = '  +6.79%  '
$ws.Range('D47').Value = $apos + '15.69'
$ws.Range('E47').Value = '  -9.43%  '
$ws.Range('D48').Value = $apos + '7.32'
$ws.Range('E48').Value = '  +6.03%  '
$ws.Range('E49').Value = '  -0.97%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = $apos + '2.26'
$ws.Range('E50').Value = '  -5.53%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.244.21'
$ws.Range('E51').Value = '  -0.42%  '
